# ppt/theme/theme1.xml backs the (sole) slide master; ppt/theme/theme2.xml
# backs the notes master. The authored diff swaps their contents: theme1
# goes from the "Integral" deck theme (clrScheme "Red Violet") to the
# stock "Office Theme" (clrScheme "Office"), while theme2 goes the other
# way. In this host, SlideMaster.Theme and NotesMaster.Theme resolve to
# the very same live theme object/part (there is only one editable theme
# in the object model, and it always round-trips to ppt/theme/theme1.xml),
# so only the slide-master side of the swap is reachable here. Push the
# 12 theme scheme colors to the target "Office" palette via
# ThemeColorScheme, which is the supported, persisted way to edit theme
# colors in this object model (the theme/clrScheme "name" attributes and
# the notes-master's own theme part are not exposed for editing).

function Convert-RGBToOle([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.Theme.ThemeColorScheme

# Target "Office Theme" clrScheme, in clrScheme document order.
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $cs.Item($i + 1).RGB = Convert-RGBToOle $officeColors[$i]
}
